$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; this shifts the existing rows 33-44
# down to 34-45 (along with their formatting), matching the weekly-update
# pattern in the diff (a brand-new record lands at row 33, the rest slide
# down one row, and the former last row, 44, becomes the new row 45).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with this week's new record.
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value = 44543
$ws.Cells.Item(33, 5).Value = 15
$ws.Cells.Item(33, 6).Value = 100112031
$ws.Cells.Item(33, 7).Value = "Poroto verde"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 1300
$ws.Cells.Item(33, 11).Value = 400
$ws.Cells.Item(33, 12).Value = 450
$ws.Cells.Item(33, 13).Value = 425
$ws.Cells.Item(33, 14).Value = "`$/kilo"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 425
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
